$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5720, 45982),
    @(5670, 45982.01041666666),
    @(5630, 45982.02083333334),
    @(5600, 45982.03125),
    @(5580, 45982.04166666666),
    @(5550, 45982.05208333334),
    @(5530, 45982.0625),
    @(5520, 45982.07291666666),
    @(5510, 45982.08333333334),
    @(5520, 45982.11458333334),
    @(5530, 45982.125),
    @(5560, 45982.13541666666),
    @(5590, 45982.14583333334),
    @(5640, 45982.15625),
    @(5700, 45982.16666666666),
    @(5780, 45982.17708333334),
    @(5880, 45982.1875),
    @(6000, 45982.19791666666),
    @(6130, 45982.20833333334),
    @(6280, 45982.21875),
    @(6440, 45982.22916666666),
    @(6610, 45982.23958333334),
    @(6770, 45982.25),
    @(6930, 45982.26041666666),
    @(7080, 45982.27083333334),
    @(7210, 45982.28125),
    @(7310, 45982.29166666666),
    @(7390, 45982.30208333334),
    @(7450, 45982.3125),
    @(7470, 45982.32291666666),
    @(7450, 45982.34375),
    @(7400, 45982.35416666666),
    @(7350, 45982.36458333334),
    @(7280, 45982.375),
    @(7220, 45982.38541666666),
    @(7170, 45982.39583333334),
    @(7120, 45982.40625),
    @(7080, 45982.41666666666),
    @(7050, 45982.42708333334),
    @(7030, 45982.4375),
    @(7010, 45982.44791666666),
    @(7000, 45982.45833333334),
    @(6990, 45982.47916666666),
    @(7000, 45982.52083333334),
    @(7010, 45982.53125),
    @(7030, 45982.54166666666),
    @(7060, 45982.55208333334),
    @(7090, 45982.5625),
    @(7120, 45982.57291666666),
    @(7160, 45982.58333333334),
    @(7200, 45982.59375),
    @(7250, 45982.60416666666),
    @(7310, 45982.61458333334),
    @(7380, 45982.625),
    @(7460, 45982.63541666666),
    @(7550, 45982.64583333334),
    @(7640, 45982.65625),
    @(7730, 45982.66666666666),
    @(7800, 45982.67708333334),
    @(7860, 45982.6875),
    @(7890, 45982.69791666666),
    @(7870, 45982.71875),
    @(7850, 45982.72916666666),
    @(7830, 45982.73958333334),
    @(7800, 45982.75),
    @(7770, 45982.76041666666),
    @(7740, 45982.77083333334),
    @(7700, 45982.78125),
    @(7660, 45982.79166666666),
    @(7600, 45982.80208333334),
    @(7530, 45982.8125),
    @(7440, 45982.82291666666),
    @(7320, 45982.83333333334),
    @(7180, 45982.84375),
    @(7060, 45982.85416666666),
    @(6940, 45982.86458333334),
    @(6780, 45982.875),
    @(6640, 45982.88541666666),
    @(6490, 45982.89583333334),
    @(6360, 45982.90625),
    @(6210, 45982.91666666666),
    @(6070, 45982.92708333334),
    @(5940, 45982.9375),
    @(5820, 45982.94791666666),
    @(5800, 45982.95833333334),
    @(5750, 45982.96875),
    @(5690, 45982.97916666666),
    @(5640, 45982.98958333334)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
